# Fix: Signature now registers for specific Plan Date, not Upload Date.
# The "fecha" (date) and ticket_id for both signature rows move forward by
# one day/id, and the "cliente" values are replaced with new placeholder
# test data (HOLA / CHAO) in place of the old ECOTRANS / AGRETRANS values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Juan Perez signature line)
$ws.Range("A2").Value = 46021
$ws.Range("B2").Value = 413345734
$ws.Range("K2").Value = "HOLA"

# Row 3 (Pedro Pascal signature line)
$ws.Range("A3").Value = 46021
$ws.Range("B3").Value = 413235734
$ws.Range("K3").Value = "CHAO"

# Move the active selection from J3 to A4
$ws.Range("A4").Select()
